# Generate Report for Handoff
# Adds two new file entries (77516efa-... and c0d11cbe-...) to the
# Overview / zh-cn / de-de worksheets, mirroring the existing
# 51c1ed14-... "Ready for handoff" rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Duplicate row 3 formatting down into rows 4 and 5 (keeps the same
# cell styles as the existing "Ready for handoff" row).
$ov.Rows(3).Copy()
$ov.Rows(4).Insert()
$ov.Rows(3).Copy()
$ov.Rows(5).Insert()

$ov.Range("A4").Value = "77516efa-27d1-4224-adc6-edb729d55a52.md"
$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = "2016-32-21 06:32:01"

$ov.Range("A5").Value = "c0d11cbe-5dd4-426c-9c34-056135eaad72.md"
$ov.Range("B5").Value = "Ready for handoff"
$ov.Range("C5").Value = "Ready for handoff"
$ov.Range("D5").Value = "2016-32-21 06:32:01"

$ov.Hyperlinks.Add($ov.Range("A4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e7eb0701391ebe53ec8a890c946bde73f15b9612/e2e/77516efa-27d1-4224-adc6-edb729d55a52.md", "", "", "77516efa-27d1-4224-adc6-edb729d55a52.md") | Out-Null
$ov.Hyperlinks.Add($ov.Range("A5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/eb58abb7e8c0147dbaa38dabadf7f0aa3d0e72eb/e2e/c0d11cbe-5dd4-426c-9c34-056135eaad72.md", "", "", "c0d11cbe-5dd4-426c-9c34-056135eaad72.md") | Out-Null

Write-Host "Overview sheet updated"

# ---------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Rows(3).Copy()
$zh.Rows(4).Insert()
$zh.Rows(3).Copy()
$zh.Rows(5).Insert()

$zh.Range("A4").Value = "77516efa-27d1-4224-adc6-edb729d55a52.md"
$zh.Range("B4").Value = ".md"
$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("D4").Value = "77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.zh-cn.xlf"
$zh.Range("E4").Value = "2016-03-21 06:31:58"
$zh.Range("H4").Value = "0001-01-01 00:00:00"
$zh.Range("I4").Value = "Include"

$zh.Range("A5").Value = "c0d11cbe-5dd4-426c-9c34-056135eaad72.md"
$zh.Range("B5").Value = ".md"
$zh.Range("C5").Value = "Ready for handoff"
$zh.Range("D5").Value = "c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.zh-cn.xlf"
$zh.Range("E5").Value = "2016-03-21 06:31:58"
$zh.Range("H5").Value = "0001-01-01 00:00:00"
$zh.Range("I5").Value = "Include"

$zh.Hyperlinks.Add($zh.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4aa7397be0f52fb2ea1fc894b587b859155bd311/e2e/77516efa-27d1-4224-adc6-edb729d55a52.md", "", "", "77516efa-27d1-4224-adc6-edb729d55a52.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/4aa7397be0f52fb2ea1fc894b587b859155bd311/e2e/77516efa-27d1-4224-adc6-edb729d55a52.md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/14059749f5738d849f44245fab459c101b6f3976/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.zh-cn.xlf", "", "", "77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.zh-cn.xlf") | Out-Null

$zh.Hyperlinks.Add($zh.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/09f54847eb7844cdee5491513883bfffc59f153c/e2e/c0d11cbe-5dd4-426c-9c34-056135eaad72.md", "", "", "c0d11cbe-5dd4-426c-9c34-056135eaad72.md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/09f54847eb7844cdee5491513883bfffc59f153c/e2e/c0d11cbe-5dd4-426c-9c34-056135eaad72.md", "", "", ".md") | Out-Null
$zh.Hyperlinks.Add($zh.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a8d1735c05efac0b374da2c9207882fc2077e875/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.zh-cn.xlf", "", "", "c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.zh-cn.xlf") | Out-Null

Write-Host "zh-cn sheet updated"

# ---------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Rows(3).Copy()
$de.Rows(4).Insert()
$de.Rows(3).Copy()
$de.Rows(5).Insert()

$de.Range("A4").Value = "77516efa-27d1-4224-adc6-edb729d55a52.md"
$de.Range("B4").Value = ".md"
$de.Range("C4").Value = "Ready for handoff"
$de.Range("D4").Value = "77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.de-de.xlf"
$de.Range("E4").Value = "2016-03-21 06:32:01"
$de.Range("H4").Value = "0001-01-01 00:00:00"
$de.Range("I4").Value = "Include"

$de.Range("A5").Value = "c0d11cbe-5dd4-426c-9c34-056135eaad72.md"
$de.Range("B5").Value = ".md"
$de.Range("C5").Value = "Ready for handoff"
$de.Range("D5").Value = "c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.de-de.xlf"
$de.Range("E5").Value = "2016-03-21 06:32:01"
$de.Range("H5").Value = "0001-01-01 00:00:00"
$de.Range("I5").Value = "Include"

$de.Hyperlinks.Add($de.Range("A4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/021bc52b530b130ce489a6497fb22b6c51cd0627/e2e/77516efa-27d1-4224-adc6-edb729d55a52.md", "", "", "77516efa-27d1-4224-adc6-edb729d55a52.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B4"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/021bc52b530b130ce489a6497fb22b6c51cd0627/e2e/77516efa-27d1-4224-adc6-edb729d55a52.md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/213075b71843fcf759236f50ef98633c5155e755/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.de-de.xlf", "", "", "77516efa-27d1-4224-adc6-edb729d55a52.670e8868800d96ac0ca6485d72ea597b44334d61.de-de.xlf") | Out-Null

$de.Hyperlinks.Add($de.Range("A5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/cf984b7e8edb5365c89389b9c326d1f2377c9063/e2e/c0d11cbe-5dd4-426c-9c34-056135eaad72.md", "", "", "c0d11cbe-5dd4-426c-9c34-056135eaad72.md") | Out-Null
$de.Hyperlinks.Add($de.Range("B5"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/cf984b7e8edb5365c89389b9c326d1f2377c9063/e2e/c0d11cbe-5dd4-426c-9c34-056135eaad72.md", "", "", ".md") | Out-Null
$de.Hyperlinks.Add($de.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b9c0adb741c00f0416d4fe0a8d0791e055231475/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.de-de.xlf", "", "", "c0d11cbe-5dd4-426c-9c34-056135eaad72.5e1f1fcd264199f25aee7d2e78d4f611949b9f08.de-de.xlf") | Out-Null

Write-Host "de-de sheet updated"
